# society import data feature improvement - area, subarea, city, state fixes
#
# - E2 "Supplier Code" value for the Jelly Beans row: KSH -> ASD
# - E3 "Supplier Code" value for the Choco Pie row:   MIT -> ZXC
# - row 3 height adjusted slightly (22.9 -> 23.05)
# - view scrolled back to top-left (A1) with the cell cursor left on E4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data fixes -----------------------------------------------------
$ws.Range("E2").Value = "ASD"
$ws.Range("E3").Value = "ZXC"

# --- row height -------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 23.05

# --- view / selection ---------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

$ws.Range("E4").Select()
